{"js": "// Remove the sentence(s) about the anchoring regression's Q1/Q2 efficacy\n// evaluation (including the two inline math expressions for Q1 and Q2)\n// that used to follow \"... individual \\u03c0 terms.\" and precede\n// \"This is justified because \" at the end of the paragraph.\n//\n// Original text (spanning several runs + two <m:oMath> fields):\n//   \" Within our anchoring regression, we can evaluate the efficacy of\n//   Q1 and Q2 in measuring substitution quality by examining the usual\n//   regression coefficients and p-values. This is justified because \"\n//\n// We locate the unique start/end text anchors with Body.search (which\n// matches plain text across run boundaries and skips over the embedded\n// math the same way Range.text linearizes it), expand a range between\n// them (this naturally swallows the two oMath fields in between), and\n// delete it in one shot.\n\nconst body = context.document.body;\n\nconst startResults = body.search(\" Within our anchoring regression\", { matchCase: false });\nconst endResults = body.search(\"This is justified because \", { matchCase: false });\nstartResults.load(\"items\");\nendResults.load(\"items\");\nawait context.sync();\n\nif (startResults.items.length === 0 || endResults.items.length === 0) {\n  throw new Error(\"Could not locate the anchoring-regression sentence to remove.\");\n}\n\nconst startRange = startResults.items[0];\nconst endRange = endResults.items[0];\n\nconst fullRange = startRange.expandTo(endRange);\nfullRange.delete();\n\nawait context.sync();\n", "ps1": "# Remove the sentence(s) about the anchoring regression's Q1/Q2 efficacy\n# evaluation (including the two inline math expressions for Q1 and Q2)\n# that used to follow \"... individual pi terms.\" and precede\n# \"This is justified because \" at the end of the paragraph.\n#\n# Original text (spanning several runs + two <m:oMath> fields):\n#   \" Within our anchoring regression, we can evaluate the efficacy of\n#   Q1 and Q2 in measuring substitution quality by examining the usual\n#   regression coefficients and p-values. This is justified because \"\n#\n# We Find.Execute a unique start anchor and a unique end anchor, capture\n# their Start/End character positions, build one contiguous Range that\n# spans both (and therefore the embedded math fields in between), and\n# delete it in one shot.\n\n$d = $word.ActiveDocument\n\n$startRange = $d.Content\n$startRange.Find.Execute(\" Within our anchoring regression\")\n$startPos = $startRange.Start\n\n$endRange = $d.Content\n$endRange.Find.Execute(\"This is justified because \")\n$endPos = $endRange.End\n\n$fullRange = $d.Range($startPos, $endPos)\n$fullRange.Delete()\n"}
